# This script reproduces the edits described by the target diff:
#  - Mark the existing "完成" status in column C for rows 94-99
#    (the six sign-off rows of the "2017.9.26" weekly block).
#  - Append a brand-new weekly block (rows 102-110) for
#    "日期：2017.9.27 第五周 周三", cloned from the formatting of the
#    preceding block (rows 92-100) so fonts / fills / borders / merges
#    all line up with the rest of the sheet.
#  - Update the sheet selection / scroll position to mirror where the
#    author ended up after typing the new block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Fill in "完成" for the six rows of the 2017.9.26 block that were
#    left blank in the "完成情况" column.
# ---------------------------------------------------------------------
$ws.Range("C94").Value = "完成"
$ws.Range("C95").Value = "完成"
$ws.Range("C96").Value = "完成"
$ws.Range("C97").Value = "完成"
$ws.Range("C98").Value = "完成"
$ws.Range("C99").Value = "完成"

# ---------------------------------------------------------------------
# 2) Clone the formatting of the previous weekly block (rows 92-100)
#    down onto the new block (rows 102-110), so the new rows pick up
#    identical fonts / alignment / fills / borders without us having to
#    hand roll style indices.
# ---------------------------------------------------------------------
$ws.Range("A92:D100").Copy()
$ws.Range("A102").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Row heights: the header/label rows keep the standard 22.5pt height,
# while the content rows grow to fit their (longer) wrapped text, just
# like the author's manual edit did.
$ws.Rows(102).RowHeight = 22.5
$ws.Rows(103).RowHeight = 22.5
$ws.Rows(104).RowHeight = 51
$ws.Rows(105).RowHeight = 77
$ws.Rows(106).RowHeight = 58
$ws.Rows(107).RowHeight = 61
$ws.Rows(108).RowHeight = 67
$ws.Rows(109).RowHeight = 48
$ws.Rows(110).RowHeight = 35

# ---------------------------------------------------------------------
# 3) New section header + merges.
# ---------------------------------------------------------------------
$ws.Range("A102").Value = "日期：2017.9.27 第五周 周三"
$ws.Range("A102:D102").Merge()

# Column header row, identical to every other block's header row.
$ws.Range("A103").Value = "人员"
$ws.Range("B103").Value = "计划任务"
$ws.Range("C103").Value = "完成情况"
$ws.Range("D103").Value = "备注"

# ---------------------------------------------------------------------
# 4) Member rows for the new block.
# ---------------------------------------------------------------------
$ws.Range("A104").Value = "李杰"
$ws.Range("B104").Value = "编写后台“用户管理”的数据交互"

$ws.Range("A105").Value = "周振朋"
$ws.Range("B105").Value = "学习HBuilder软件操作并学习h5教程，完成视频播放demo"

$ws.Range("A106").Value = "禤锦辉"
$ws.Range("B106").Value = "学习HBuilder软件操作并学习h5教程"

$ws.Range("A107").Value = "柯新钿"
$ws.Range("B107").Value = "学习HBuilder软件操作并学习h5教程，完成拖放demo"

$ws.Range("A108").Value = "冯文雄"
$ws.Range("B108").Value = "编写后台“商品类别管理”的数据交互"

$ws.Range("A109").Value = "阿卜力孜"
$ws.Range("B109").Value = "学习HBuilder软件操作并学习h5教程"

# Summary / sign-off row for the new block.
$ws.Range("A110").Value = "总结："
$ws.Range("A110:D110").Merge()

# ---------------------------------------------------------------------
# 5) Leave the sheet scrolled/selected where the author ended up after
#    entering the new block.
# ---------------------------------------------------------------------
$ws.Range("A101").Select()
$excel.ActiveWindow.ScrollRow = 101
$ws.Range("A110:D110").Select()

Write-Output "edit complete"
